$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp value in A19 (tiny precision change)
$ws.Cells.Item(19, 1).Value = 45876.70856442129

# Add new row 20 with data
$ws.Cells.Item(20, 1).Value = 45876.75020091004
$ws.Cells.Item(20, 1).NumberFormat = $ws.Cells.Item(19, 1).NumberFormat

$ws.Cells.Item(20, 2).Value = 2025
$ws.Cells.Item(20, 3).Value = 28
$ws.Cells.Item(20, 4).Value = 16.67
$ws.Cells.Item(20, 5).Value = 84.15000000000001
$ws.Cells.Item(20, 6).Value = 15.82
$ws.Cells.Item(20, 7).Value = 10.39
$ws.Cells.Item(20, 8).Value = "ESE"
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = "18:00:17"
